# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" > "Impact" bullet list from
# job-duty style statements into impact-focused accomplishment statements.
#
# The six original bullet paragraphs become four new bullet paragraphs:
#   1. "Built real-time FEC analysis systems..."        -> "Algorithmic innovation: ..."
#   2. "Built cloud-based data warehouse solutions..."   -> "$4.7M savings enabled nonprofit access"
#   3. "Designed ETL pipelines..."                       -> (removed)
#   4. "Trigonometric algorithm for boundary estimation..." -> "Breakthrough demographic discovery: ..."
#   5. "Built redistricting platform..."                 -> (removed)
#   6. "Discovered systematic race coding errors..."     -> "178% accuracy improvement in racial classification algorithms"

$d = $word.ActiveDocument

function Set-ParagraphText($paragraph, $newText) {
    $r = $paragraph.Range
    $rng = $d.Range($r.Start, $r.End - 1)
    $rng.Text = $newText
}

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading, then the "Impact" sub-heading
# that immediately follows it, so the right occurrence of each (duplicated)
# bullet string is targeted regardless of paragraph numbering elsewhere.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# "Impact" sub-heading directly follows; bullets start right after that.
$firstBulletIndex = $headingIndex + 2

$p1 = $d.Paragraphs.Item($firstBulletIndex)
$p2 = $d.Paragraphs.Item($firstBulletIndex + 1)
$p3 = $d.Paragraphs.Item($firstBulletIndex + 2)
$p4 = $d.Paragraphs.Item($firstBulletIndex + 3)
$p5 = $d.Paragraphs.Item($firstBulletIndex + 4)
$p6 = $d.Paragraphs.Item($firstBulletIndex + 5)

Write-Host "Bullet 1:" $p1.Range.Text
Write-Host "Bullet 2:" $p2.Range.Text
Write-Host "Bullet 3:" $p3.Range.Text
Write-Host "Bullet 4:" $p4.Range.Text
Write-Host "Bullet 5:" $p5.Range.Text
Write-Host "Bullet 6:" $p6.Range.Text

# Sanity-check the existing content before mutating, so the script fails
# loudly instead of silently editing the wrong paragraphs.
if ($p3.Range.Text -notmatch "Designed ETL pipelines") {
    throw "Unexpected text in bullet 3 (expected the ETL pipelines bullet)"
}
if ($p5.Range.Text -notmatch "Built redistricting platform") {
    throw "Unexpected text in bullet 5 (expected the redistricting platform bullet)"
}

# Update the four bullets that survive, then delete the two that don't -
# deleting from the bottom up so earlier paragraph references stay valid.
$bullet = [char]0x2022
Set-ParagraphText $p1 "$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
Set-ParagraphText $p2 "$bullet `$4.7M savings enabled nonprofit access"
Set-ParagraphText $p4 "$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
Set-ParagraphText $p6 "$bullet 178% accuracy improvement in racial classification algorithms"

$p5.Range.Delete()
$p3.Range.Delete()

Write-Host "Result bullets:"
for ($i = $firstBulletIndex; $i -le $firstBulletIndex + 3; $i++) {
    Write-Host "  " $d.Paragraphs.Item($i).Range.Text
}
